# Applies the LinuxForHealth rebrand + version bump edit to the
# StructureDefinition-quality-measures workbook.

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet ----
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/quality-measures"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---- Elements sheet ----
$elem = $wb.Worksheets.Item("Elements")

# The "Extension" row no longer repeats the ele-1/ext-1 constraint text
# (it now lives only on the "Extension.extension" row).
$elem.Range("AI2").Value = ""

# The canonical URL referenced from the Extension.url fixed value and the
# Extension.value[x] type both move from ibm.com to linuxforhealth.org.
$elem.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/quality-measures"
$elem.Range("J6").Value = "Reference(http://linuxforhealth.org/fhir/cdm/StructureDefinition/quality-measure-list)
"

$wb.Save()
